# Generate Report for Handoff
# Adds a new row (row 3) to each of the three worksheets (Overview, zh-cn, de-de)
# describing the handoff of f30d5173-3193-4915-b48c-f81210d73ab0, mirroring the
# existing row for 560101b1-a1ef-4878-b250-f85ebe891b31.

$wb = $excel.ActiveWorkbook

$oldGuid = "560101b1-a1ef-4878-b250-f85ebe891b31"
$newGuid = "f30d5173-3193-4915-b48c-f81210d73ab0"
$newHash = "23c9c913752ae51596a2004e31e078ee7e3796c1"

$hyperlinkColor = 15570276   # BGR long for RGB FF6495ED (matches the workbook's "HyperLink" style)

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A3").Value = "$newGuid.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-34-20 16:34:11"

$lnk = $ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/22f74c5e4f533a8f0060f2cda23ef5660869447d/e2e/$newGuid.md", "", "", "$newGuid.md")
$ws1.Range("A3").Font.Underline = 2
$ws1.Range("A3").Font.Color = $hyperlinkColor

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A3").Value = "$newGuid.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "$newGuid.$newHash.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-20 16:34:08"
$ws2.Range("H3").Value = "0001-01-01 00:00:00"
$ws2.Range("I3").Value = "Include"

$lnk = $ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/22f74c5e4f533a8f0060f2cda23ef5660869447d/e2e/$newGuid.md", "", "", "$newGuid.md")
$ws2.Range("A3").Font.Underline = 2
$ws2.Range("A3").Font.Color = $hyperlinkColor

$lnk = $ws2.Hyperlinks.Add($ws2.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/22f74c5e4f533a8f0060f2cda23ef5660869447d/e2e/$newGuid.md", "", "", ".md")
$ws2.Range("B3").Font.Underline = 2
$ws2.Range("B3").Font.Color = $hyperlinkColor

$lnk = $ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c162d1a28e6467a7c4e80d32583a4e5e14fc0adb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newGuid.$newHash.zh-cn.xlf", "", "", "$newGuid.$newHash.zh-cn.xlf")
$ws2.Range("D3").Font.Underline = 2
$ws2.Range("D3").Font.Color = $hyperlinkColor

$ws2.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A3").Value = "$newGuid.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "$newGuid.$newHash.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-20 16:34:11"
$ws3.Range("H3").Value = "0001-01-01 00:00:00"
$ws3.Range("I3").Value = "Include"

$lnk = $ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/22f74c5e4f533a8f0060f2cda23ef5660869447d/e2e/$newGuid.md", "", "", "$newGuid.md")
$ws3.Range("A3").Font.Underline = 2
$ws3.Range("A3").Font.Color = $hyperlinkColor

$lnk = $ws3.Hyperlinks.Add($ws3.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/22f74c5e4f533a8f0060f2cda23ef5660869447d/e2e/$newGuid.md", "", "", ".md")
$ws3.Range("B3").Font.Underline = 2
$ws3.Range("B3").Font.Color = $hyperlinkColor

$lnk = $ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9abb8d041e474ecd509046179e256280413fabe3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newGuid.$newHash.de-de.xlf", "", "", "$newGuid.$newHash.de-de.xlf")
$ws3.Range("D3").Font.Underline = 2
$ws3.Range("D3").Font.Color = $hyperlinkColor

$ws3.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
